$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells are stored as plain text in the source data so that
# trailing zeros (e.g. "54.60") are preserved. Pre-set NumberFormat to Text
# ("@") on exactly the Price cells we are about to rewrite with a plain
# decimal value, so Excel does not silently convert them to numbers and lose
# the trailing zero / exact formatting. Cells whose new value still contains
# two dots (e.g. "41.739.90") are left alone - Excel already treats those as
# text because they cannot be parsed as a number.
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D9:D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21:D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31:D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46:D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "41.739.90"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "2.176.00"
$ws.Range("E3").Value = "  -2.66%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "238.63"
$ws.Range("E5").Value = "  -1.72%  "

$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("D7").Value = "72.69"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -2.80%  "

$ws.Range("D10").Value = "40.53"
$ws.Range("E10").Value = "  -4.15%  "

$ws.Range("D11").Value = "0.0912"
$ws.Range("E11").Value = "  -4.75%  "

$ws.Range("D12").Value = "54.60"
$ws.Range("E12").Value = "  -3.45%  "

$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("E14").Value = "  -3.05%  "

$ws.Range("D15").Value = "2.503.60"
$ws.Range("E15").Value = "  -2.64%  "

$ws.Range("D16").Value = "14.39"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "2.169.90"
$ws.Range("E17").Value = "  -2.68%  "

$ws.Range("D18").Value = "0.784"
$ws.Range("E18").Value = "  -6.40%  "

$ws.Range("D19").Value = "41.662.27"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").Value = "70.22"
$ws.Range("E21").Value = "  -3.54%  "

$ws.Range("D22").Value = "5.82"
$ws.Range("E22").Value = "  -6.59%  "

$ws.Range("D23").Value = "10.11"
$ws.Range("E23").Value = "  -12.22%  "

$ws.Range("D24").Value = "226.65"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  -4.98%  "

$ws.Range("E28").Value = "  -9.97%  "

$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -3.41%  "

$ws.Range("E30").Value = "  -1.10%  "

$ws.Range("D31").Value = "170.75"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("D32").Value = "19.89"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("D33").Value = "32.78"
$ws.Range("E33").Value = "  +9.69%  "

$ws.Range("D34").Value = "0.0778"
$ws.Range("E34").Value = "  -3.29%  "

$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -5.63%  "

$ws.Range("D36").Value = "0.121"
$ws.Range("E36").Value = "  -3.35%  "

$ws.Range("D37").Value = "4.33"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("E38").Value = "  -6.43%  "

$ws.Range("D39").Value = "0.0311"
$ws.Range("E39").Value = "  +2.33%  "

$ws.Range("E40").Value = "  -7.95%  "

$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("E42").Value = "  -5.77%  "

$ws.Range("D43").Value = "59.44"
$ws.Range("E43").Value = "  -8.09%  "

$ws.Range("E44").Value = "  -4.55%  "

$ws.Range("E45").Value = "  -2.92%  "

$ws.Range("D46").Value = "0.0965"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").Value = "97.56"
$ws.Range("E47").Value = "  -6.64%  "

$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("E49").Value = "  -4.64%  "

$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -5.94%  "

$ws.Range("E51").Value = "  -2.04%  "
